$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "/flashcard.html"

$ws.Range("I7").Select()
